$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.272216320037842
$ws.Range("B1").Value = 2.682085037231445
$ws.Range("C1").Value = 2.836045026779175
$ws.Range("D1").Value = 3.551714420318604
$ws.Range("E1").Value = 1.617812275886536
